$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 182, shifting the existing
# rows 182-189 down to 184-191 (matching the target dimension A1:R191).
$ws.Rows("182:183").Insert()

# New row 182: Comercializadora del Agro de Limari - Espanola / Primera, 2022-05-25
$ws.Cells.Item(182, 1).Value = 2
$ws.Cells.Item(182, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(182, 3).Value = "Coquimbo"
$ws.Cells.Item(182, 4).Value = 44706
$ws.Cells.Item(182, 5).Value = 4
$ws.Cells.Item(182, 6).Value = 100112013
$ws.Cells.Item(182, 7).Value = "Alcachofa"
$ws.Cells.Item(182, 8).Value = "Española"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 500
$ws.Cells.Item(182, 11).Value = 16000
$ws.Cells.Item(182, 12).Value = 17000
$ws.Cells.Item(182, 13).Value = 16500
$ws.Cells.Item(182, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(182, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(182, 16).Value = 550
$ws.Cells.Item(182, 17).Value = 30
$ws.Cells.Item(182, 18).Value = "Hortaliza"

# New row 183: Comercializadora del Agro de Limari - Espanola / Segunda, 2022-05-25
$ws.Cells.Item(183, 1).Value = 2
$ws.Cells.Item(183, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(183, 3).Value = "Coquimbo"
$ws.Cells.Item(183, 4).Value = 44706
$ws.Cells.Item(183, 5).Value = 4
$ws.Cells.Item(183, 6).Value = 100112013
$ws.Cells.Item(183, 7).Value = "Alcachofa"
$ws.Cells.Item(183, 8).Value = "Española"
$ws.Cells.Item(183, 9).Value = "Segunda"
$ws.Cells.Item(183, 10).Value = 300
$ws.Cells.Item(183, 11).Value = 14000
$ws.Cells.Item(183, 12).Value = 15000
$ws.Cells.Item(183, 13).Value = 14500
$ws.Cells.Item(183, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(183, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(183, 16).Value = 362
$ws.Cells.Item(183, 17).Value = 40
$ws.Cells.Item(183, 18).Value = "Hortaliza"
